$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row heights
$ws.Rows.Item(1).RowHeight = 17
$ws.Range("A2:A22").EntireRow.RowHeight = 30
$ws.Rows.Item(23).RowHeight = 17

# Update the view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 15
$ws.Range("A1:C23").Select()
$ws.Range("C23").Activate()
